$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the truck insurance smoke test entries: drop "Automobile" from the name
$ws.Range("A4").Value = "103_TruckInsurance_001_SmokeTest"
$ws.Range("B4").Value = "var103_TruckInsurance_001_SmokeTest"
$ws.Range("E4").Value = "103_TruckInsurance_001_SmokeTest"

# Update the selected cell / active cell shown when the file is saved
$ws.Range("A6").Select()
